# Update recomputed TPM-derived NATMI ligand-receptor statistics
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G"=0.8999933333333333; "H"=2.69998; "I"=0.01781586806325543; "J"=0.01781586806325543; "M"=118.0346986666667; "N"=354.104096; "O"=0.2666057129183408; "P"=0.2666057129183408; "Q"=106.2304419020089; "R"=956.0739771180799; "S"=0.004749812206263314; "T"=0.004749812206263315 }
    3 = @{ "G"=0.8999933333333333; "H"=2.69998; "I"=0.01781586806325543; "J"=0.01781586806325543; "O"=0.4881754016778185; "P"=0.4881754016778186; "Q"=194.51604422974; "R"=1750.64439806766; "S"=0.008697268548018737; "T"=0.008697268548018739 }
    4 = @{ "G"=0.8999933333333333; "H"=2.69998; "I"=0.01781586806325543; "J"=0.01781586806325543; "M"=45.876452; "N"=137.629356; "O"=0.1036214293744632; "P"=0.1036214293744632; "Q"=41.28850095698667; "R"=371.59650861288; "S"=0.001846105714261376; "T"=0.001846105714261377 }
    5 = @{ "G"=0.8999933333333333; "H"=2.69998; "I"=0.01781586806325543; "J"=0.01781586806325543; "M"=62.68962833333333; "N"=188.068885; "O"=0.1415974560293775; "P"=0.1415974560293775; "Q"=56.42024756914444; "R"=507.7822281223; "S"=0.002522681594712; "T"=0.002522681594712002 }
    6 = @{ "I"=0.7949938412397365; "J"=0.7949938412397366; "M"=118.0346986666667; "N"=354.104096; "O"=0.2666057129183408; "P"=0.2666057129183408; "Q"=4740.299308707442; "R"=42662.69377836698; "S"=0.2119498998094102; "T"=0.2119498998094103 }
    7 = @{ "I"=0.7949938412397365; "J"=0.7949938412397366; "O"=0.4881754016778185; "P"=0.4881754016778186; "S"=0.3880964377786003; "T"=0.3880964377786004 }
    8 = @{ "I"=0.7949938412397365; "J"=0.7949938412397366; "M"=45.876452; "N"=137.629356; "O"=0.1036214293744632; "P"=0.1036214293744632; "Q"=1842.408343970838; "R"=16581.67509573754; "S"=0.08237839817315655; "T"=0.08237839817315658 }
    9 = @{ "I"=0.7949938412397365; "J"=0.7949938412397366; "M"=62.68962833333333; "N"=188.068885; "O"=0.1415974560293775; "P"=0.1415974560293775; "Q"=2517.629181998729; "R"=22658.66263798856; "S"=0.1125691054785695; "T"=0.1125691054785695 }
    10 = @{ "G"=5.293300666666666; "H"=15.879902; "I"=0.104783827617029; "J"=0.104783827617029; "M"=118.0346986666667; "N"=354.104096; "O"=0.2666057129183408; "P"=0.2666057129183408; "Q"=624.7931491420657; "R"=5623.138342278591; "S"=0.02793596706415055; "T"=0.02793596706415055 }
    11 = @{ "G"=5.293300666666666; "H"=15.879902; "I"=0.104783827617029; "J"=0.104783827617029; "O"=0.4881754016778185; "P"=0.4881754016778186; "Q"=1144.043926175726; "R"=10296.39533558153; "S"=0.05115288713628243; "T"=0.05115288713628244 }
    12 = @{ "G"=5.293300666666666; "H"=15.879902; "I"=0.104783827617029; "J"=0.104783827617029; "M"=45.876452; "N"=137.629356; "O"=0.1036214293744632; "P"=0.1036214293744632; "Q"=242.8378539559013; "R"=2185.540685603112; "S"=0.0108578499930039; "T"=0.0108578499930039 }
    13 = @{ "G"=5.293300666666666; "H"=15.879902; "I"=0.104783827617029; "J"=0.104783827617029; "M"=62.68962833333333; "N"=188.068885; "O"=0.1415974560293775; "P"=0.1415974560293775; "Q"=331.8350514499189; "R"=2986.515463049269; "S"=0.01483712342359213; "T"=0.01483712342359213 }
    14 = @{ "G"=4.162877000000001; "H"=12.488631; "I"=0.08240646307997901; "J"=0.08240646307997902; "M"=118.0346986666667; "N"=354.104096; "O"=0.2666057129183408; "P"=0.2666057129183408; "Q"=491.3639322813974; "R"=4422.275390532576; "S"=0.02197003383851674; "T"=0.02197003383851674 }
    15 = @{ "G"=4.162877000000001; "H"=12.488631; "I"=0.08240646307997901; "J"=0.08240646307997902; "O"=0.4881754016778185; "P"=0.4881754016778186; "Q"=899.7248498007032; "R"=8097.523648206328; "S"=0.04022880821491708; "T"=0.04022880821491709 }
    16 = @{ "G"=4.162877000000001; "H"=12.488631; "I"=0.08240646307997901; "J"=0.08240646307997902; "M"=45.876452; "N"=137.629356; "O"=0.1036214293744632; "P"=0.1036214293744632; "Q"=190.978026872404; "R"=1718.802241851636; "S"=0.008539075494041353; "T"=0.008539075494041355 }
    17 = @{ "G"=4.162877000000001; "H"=12.488631; "I"=0.08240646307997901; "J"=0.08240646307997902; "M"=62.68962833333333; "N"=188.068885; "O"=0.1415974560293775; "P"=0.1415974560293775; "Q"=260.9692119273817; "R"=2348.722907346435; "S"=0.01166854553250384; "T"=0.01166854553250385 }
}

foreach ($rowNum in $updates.Keys) {
    $rowVals = $updates[$rowNum]
    foreach ($colLetter in $rowVals.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowVals[$colLetter]
    }
}

